$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "haresh"
$ws.Range("B1").Value = "parab"
